$wb = $excel.ActiveWorkbook

# Sheet 1: "24_HRS" -> rows 2..25, column A, hourly timestamps
# starting at 2023-05-06 14:03 and going back 1 hour per row.
$ws1 = $wb.Worksheets.Item("24_HRS")
$start1 = Get-Date -Year 2023 -Month 5 -Day 6 -Hour 14 -Minute 3 -Second 0
for ($i = 0; $i -lt 24; $i++) {
    $dt = $start1.AddHours(-$i)
    $ws1.Cells.Item($i + 2, 1).Value = $dt.ToString("yyyy-MM-dd HH:mm")
}

# Sheet 2: "1d_bef" -> rows 2..25, column A, hourly timestamps
# starting at 2023-05-05 15:03 and going back 1 hour per row.
# Rows 3..25 are newly populated; column B stays an empty cell like row 2,
# so copy the already-empty B2 cell down to materialize empty B cells.
$ws2 = $wb.Worksheets.Item("1d_bef")
$start2 = Get-Date -Year 2023 -Month 5 -Day 5 -Hour 15 -Minute 3 -Second 0
for ($i = 0; $i -lt 24; $i++) {
    $dt = $start2.AddHours(-$i)
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $dt.ToString("yyyy-MM-dd HH:mm")
}
$ws2.Range("B2").Copy($ws2.Range("B3:B25"))

# Sheet 3: "7d_bef" -> rows 2..25, column A, hourly timestamps
# starting at 2023-04-29 15:03 and going back 1 hour per row.
$ws3 = $wb.Worksheets.Item("7d_bef")
$start3 = Get-Date -Year 2023 -Month 4 -Day 29 -Hour 15 -Minute 3 -Second 0
for ($i = 0; $i -lt 24; $i++) {
    $dt = $start3.AddHours(-$i)
    $ws3.Cells.Item($i + 2, 1).Value = $dt.ToString("yyyy-MM-dd HH:mm")
}
